# Locate the shape that holds the "Δt 时间内所经过的距离 L 可以近似为 ..." sentence
# and merge the split runs back together, per the target diff:
#   "时间内所经过" + "的距离 "      -> "时间内所经过的距离 "
#   "可以近似为"   + " "          -> "可以近似为 "
$p = $ppt.ActivePresentation

$marker = "时间内所经过"
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                $txt = $shape.TextFrame.TextRange.Text
                if ($txt.IndexOf($marker) -ge 0) {
                    $targetSlide = $slide
                    $targetShape = $shape
                }
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$full = $tr.Text

# --- Edit 2 first (further to the right) so the offsets used for edit 1 stay valid ---
# "可以近似为" (5 chars) immediately followed by a lone " " run (1 char) -> "可以近似为 "
$idx2 = $full.IndexOf("可以近似为")
$start2 = $idx2 + 1
$c2 = $tr.Characters($start2, 6)
$c2.Text = "可以近似为 "

# --- Edit 1: "时间内所经过" (6 chars) followed by "的距离 " (4 chars) -> "时间内所经过的距离 " ---
$idx1 = $full.IndexOf($marker)
$start1 = $idx1 + 1
$c1 = $tr.Characters($start1, 10)
$c1.Text = "时间内所经过的距离 "
